$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally listed 5 exposure sites (Abbotsford, Albert Park,
# Hampton, and two Springvale rows). The update trims this down to a
# single Springvale site (IKEA Springvale), keeping one current ("new")
# entry and one superseded ("old") entry for it.

# Delete the second Springvale row (old row 6, the shopping-centre
# entry) first so the row numbers above it are unaffected.
$ws.Rows("6:6").Delete()

# Delete Abbotsford (2), Albert Park (3) and Hampton (4). What remains
# is the header row plus the IKEA Springvale row (old row 5), which
# becomes the new row 2.
$ws.Rows("2:4").Delete()

# Shorten the surviving (current) Springvale row's notes.
$ws.Range("D2").Value2 = "Case shopped at store"

# Add the superseded entry for the same Springvale site, carrying the
# original, longer note and marked as "old".
$ws.Range("A3").Value2 = "Springvale"
$ws.Range("B3").Value2 = "IKEA Springvale, 917 Princes Hwy"
$ws.Range("C3").Value2 = "30/12/20, 4:00pm-6:30pm"
$ws.Range("D3").Value2 = "Case shopped at store and dined at cafe"
$ws.Range("E3").Value2 = "old"

# Re-fit the first three columns to their new (shorter) content and
# match the active selection left behind by the edit.
$ws.Columns("A").ColumnWidth = 8.86328125
$ws.Columns("B").ColumnWidth = 27.06640625
$ws.Columns("C").ColumnWidth = 22.06640625

[void]$ws.Range("D3").Select()
